$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3679053
$ws.Range("I69").Value = 14705882
$ws.Range("J69").Value = 3443.3333
$ws.Range("K69").Value = 44117646
$ws.Range("L69").Value = 10329.9999
$ws.Range("M69").Value = -44116772
$ws.Range("N69").Value = -12077.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 3679053
$ws.Range("I72").Value = 14705882
$ws.Range("J72").Value = 3443.3333
$ws.Range("K72").Value = 132352938
$ws.Range("L72").Value = 30989.9997
$ws.Range("M72").Value = -132348570
$ws.Range("N72").Value = -39725.9997

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 2781430.2
$ws.Range("I74").Value = 3033951.5
$ws.Range("J74").Value = 3696.6667
$ws.Range("K74").Value = 3033951.5
$ws.Range("L74").Value = 3696.6667
$ws.Range("M74").Value = -3033015.5
$ws.Range("N74").Value = -5568.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 2781430.2
$ws.Range("I77").Value = 3033951.5
$ws.Range("J77").Value = 3696.6667
$ws.Range("K77").Value = 15169757.5
$ws.Range("L77").Value = 18483.3335
$ws.Range("M77").Value = -15165077.5
$ws.Range("N77").Value = -27843.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1927.2927
$ws.Range("I129").Value = 887.4
$ws.Range("J129").Value = 2071.7222
$ws.Range("K129").Value = 2662.2
$ws.Range("L129").Value = 6215.1666
$ws.Range("M129").Value = 2337.8
$ws.Range("N129").Value = -16215.1666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H130").Value = 30599.6
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 30599.6
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 30599.6
$ws.Range("N130").Value = -40639.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2255.75
$ws.Range("I45").Value = 1912.6471
$ws.Range("J45").Value = 4200
$ws.Range("K45").Value = 1912.6471
$ws.Range("L45").Value = 4200
$ws.Range("M45").Value = -1535.6471
$ws.Range("N45").Value = -4954

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2568.182
$ws.Range("I61").Value = 2333.1538
$ws.Range("J61").Value = 3441.1428
$ws.Range("K61").Value = 2333.1538
$ws.Range("L61").Value = 3441.1428
$ws.Range("M61").Value = -2121.1538
$ws.Range("N61").Value = -3865.1428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2599.375
$ws.Range("I63").Value = 2599.375
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2599.375
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1913.375
$ws.Range("N63").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2599.375
$ws.Range("I66").Value = 2599.375
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 12996.875
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -9564.875
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 51881.95
$ws.Range("I74").Value = 64755.062
$ws.Range("J74").Value = 389.5
$ws.Range("K74").Value = 64755.062
$ws.Range("L74").Value = 389.5
$ws.Range("M74").Value = -63881.062
$ws.Range("N74").Value = -2137.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 51881.95
$ws.Range("I77").Value = 64755.062
$ws.Range("J77").Value = 389.5
$ws.Range("K77").Value = 323775.31
$ws.Range("L77").Value = 1947.5
$ws.Range("M77").Value = -319407.31
$ws.Range("N77").Value = -10683.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H111").Value = 37312
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 37312
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 37312
$ws.Range("N111").Value = -45492

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2568.182
$ws.Range("I136").Value = 2333.1538
$ws.Range("J136").Value = 3441.1428
$ws.Range("K136").Value = 6999.4614
$ws.Range("L136").Value = 10323.4284
$ws.Range("M136").Value = -4449.4614
$ws.Range("N136").Value = -15423.4284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10072.487
$ws.Range("I31").Value = 11084.95
$ws.Range("J31").Value = 9006.736999999999
$ws.Range("K31").Value = 11084.95
$ws.Range("L31").Value = 9006.736999999999
$ws.Range("M31").Value = -10789.95

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 10072.487
$ws.Range("I34").Value = 11084.95
$ws.Range("J34").Value = 9006.736999999999
$ws.Range("K34").Value = 11084.95
$ws.Range("L34").Value = 9006.736999999999
$ws.Range("M34").Value = -10882.95

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3637.524
$ws.Range("I58").Value = 1035.9
$ws.Range("J58").Value = 6002.636
$ws.Range("K58").Value = 1035.9
$ws.Range("L58").Value = 6002.636
$ws.Range("M58").Value = -832.9000000000001
$ws.Range("N58").Value = -6408.636

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 21740650
$ws.Range("I134").Value = 1420.375
$ws.Range("J134").Value = 71430320
$ws.Range("K134").Value = 4261.125
$ws.Range("L134").Value = 214290960
$ws.Range("M134").Value = -1726.125
$ws.Range("N134").Value = -214296030

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3637.524
$ws.Range("I136").Value = 1035.9
$ws.Range("J136").Value = 6002.636
$ws.Range("K136").Value = 3107.7
$ws.Range("L136").Value = 18007.908
$ws.Range("M136").Value = -557.7000000000003
$ws.Range("N136").Value = -23107.908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 13115944
$ws.Range("I11").Value = 18156692
$ws.Range("J11").Value = 10000
$ws.Range("K11").Value = 18156692
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = -18156553
$ws.Range("N11").Value = -10278

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 50670.668
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 50670.668
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 50670.668
$ws.Range("N18").Value = -51256.668

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4258.095
$ws.Range("I70").Value = 4226.6665
$ws.Range("J70").Value = 4300
$ws.Range("K70").Value = 4226.6665
$ws.Range("L70").Value = 4300
$ws.Range("M70").Value = -3956.6665
$ws.Range("N70").Value = -4840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4258.095
$ws.Range("I73").Value = 4226.6665
$ws.Range("J73").Value = 4300
$ws.Range("K73").Value = 4226.6665
$ws.Range("L73").Value = 4300
$ws.Range("M73").Value = -3290.6665
$ws.Range("N73").Value = -6172

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 34802.4
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 34802.4
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 34802.4
$ws.Range("N20").Value = -35254.4
$ws.Range("M20").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1828.8572
$ws.Range("I81").Value = 1666.6666
$ws.Range("J81").Value = 1873.091
$ws.Range("K81").Value = 3333.3332
$ws.Range("L81").Value = 3746.182
$ws.Range("M81").Value = -2272.3332
$ws.Range("N81").Value = -5868.182

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1828.8572
$ws.Range("I84").Value = 1666.6666
$ws.Range("J84").Value = 1873.091
$ws.Range("K84").Value = 16666.666
$ws.Range("L84").Value = 18730.91
$ws.Range("M84").Value = -11362.666
$ws.Range("N84").Value = -29338.91
